$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember what's currently in the summary block before we shift it down ---
$label46 = $ws.Range("E46").Value2
$label47 = $ws.Range("E47").Value2
$label48 = $ws.Range("E48").Value2

# --- Push the blank separator row + the three summary rows down by one row,
#     copying formatting first (bottom-up, so we never clobber a source
#     before it has been copied from) ---
$ws.Range("E48").Copy()
$ws.Range("E49").PasteSpecial(-4122) | Out-Null
$ws.Range("F48").Copy()
$ws.Range("F49").PasteSpecial(-4122) | Out-Null

$ws.Range("E47").Copy()
$ws.Range("E48").PasteSpecial(-4122) | Out-Null
$ws.Range("F47").Copy()
$ws.Range("F48").PasteSpecial(-4122) | Out-Null

$ws.Range("E46").Copy()
$ws.Range("E47").PasteSpecial(-4122) | Out-Null
$ws.Range("F46").Copy()
$ws.Range("F47").PasteSpecial(-4122) | Out-Null

$ws.Range("D45").Copy()
$ws.Range("D46").PasteSpecial(-4122) | Out-Null
$ws.Range("E45").Copy()
$ws.Range("E46").PasteSpecial(-4122) | Out-Null
$ws.Range("F45").Copy()
$ws.Range("F46").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Re-home the summary formulas/labels one row lower ---
$ws.Range("E49").Value = $label48
$ws.Range("F49").Formula = "=F48/38.5"

$ws.Range("E48").Value = $label47
$ws.Range("F48").Formula = "=F47/60"

$ws.Range("E47").Value = $label46
$ws.Range("F47").Formula = "=SUM(F2:F46)"

# Row 46 becomes the new blank separator row (formats already copied above).
$ws.Range("D46").ClearContents()
$ws.Range("E46").ClearContents()
$ws.Range("F46").ClearContents()

# --- Update the existing entry on row 44 (its end time changed) ---
$ws.Range("E44").Value = 0.54166666666666663

# --- Row 45 becomes a brand-new time entry ---
$ws.Range("A45").Value = 2014
$ws.Range("B45").Value = 3
$ws.Range("C45").Value = 4

$ws.Range("D44").Copy()
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").Copy()
$ws.Range("E45").PasteSpecial(-4122) | Out-Null
$ws.Range("F44").Copy()
$ws.Range("F45").PasteSpecial(-4122) | Out-Null
$ws.Range("G44").Copy()
$ws.Range("G45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D45").Value = 0.58333333333333337
$ws.Range("E45").Value = 0.60416666666666663
$ws.Range("F45").Formula = "=(E45-D45)*24*60"
$ws.Range("G45").Formula = "=F45/60"

$ws.Range("F45").Select()
